# PROC-106 | refact: changes simple procuration docs
#
# Removes the three paragraphs that followed the "_proc_today_." line
# (the underscore signature line + "_proc_full_name_" line + the blank
# paragraph right after it), collapsing straight into the following
# blank paragraph that precedes the "(Rep. Legal)" signature block.

$d = $word.ActiveDocument

$paras = $d.Paragraphs

# Locate the paragraph that contains the "_proc_today_." marker.
$todayIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*_proc_today_.*") {
        $todayIndex = $i
        break
    }
}

if ($todayIndex -gt 0) {
    # The three paragraphs to remove are the ones immediately following it:
    #   - the underline "_____" signature rule
    #   - the "_proc_full_name_" line
    #   - the blank paragraph right after it
    $firstToRemove = $paras.Item($todayIndex + 1)
    $lastToRemove = $paras.Item($todayIndex + 3)

    $rangeStart = $firstToRemove.Range.Start
    $rangeEnd = $lastToRemove.Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
